# Update the "Pitch correction" row (row 25) values from 15.8 to 15
# for all four datasets (columns B-E). All dependent formula cells
# (B29:E29, B30:E30, B33:E33, B34:E34, I12:L12, N12) recalculate
# automatically from these inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25:E25").Value = 15

# Restore the active cell selection recorded in the saved workbook.
$ws.Range("F25").Select()

$excel.CalculateFullRebuild()
